# Applies the "additional scraping" changes described in the commit:
#  1. Insert a new "Player Info" sheet as the first sheet in the workbook,
#     containing the player's ID/NAME/BATTING_HAND/BOWL_STYLE.
#  2. Rename the "MATCH_CARD_LINK" column to "MATCH_CODE" on both the
#     "ODI Batting" and "ODI Bowling" sheets, and replace the full
#     scorecard URL values with just the bare numeric match code.

$wb = $excel.ActiveWorkbook

# Helper: write a value that looks numeric (e.g. "4625") but must be kept
# as plain text, without leaving the cell tagged with a new/different
# style than it started with.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$battingSheet = $wb.Worksheets.Item("ODI Batting")

# --- 1. Add the new "Player Info" worksheet, placed before "ODI Batting" ---
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

# Re-fetch sheet references by name now that the sheet collection changed,
# since sheet references track position rather than identity in this
# automation environment.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Copy the existing header style (bold, centered, bordered) from the
# "ODI Batting" header row onto the new sheet's header row instead of
# building a brand-new style entry.
$battingSheet.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

Set-TextValue $playerInfo.Range("A2") "4625"
$playerInfo.Range("B2").Value = "Duanne Olivier"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium Fast"

# --- 2. Update "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$battingSheet.Range("D1").Value = "MATCH_CODE"
Set-TextValue $battingSheet.Range("D2") "4237"
Set-TextValue $battingSheet.Range("D3") "4238"

# --- 3. Update "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
Set-TextValue $bowlingSheet.Range("B2") "4237"
Set-TextValue $bowlingSheet.Range("B3") "4238"
